$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed from the missing-data sample
# (delete the higher-numbered row first so the other row index stays valid)
$ws.Rows.Item(28).Delete()   # "SC 92"
$ws.Rows.Item(26).Delete()   # "RM 232"

# Re-roll which cells in column D ("C") are treated as missing for this seed
$ws.Range("D3").Value = -14.2    # RM 8: was missing, now has a value
$ws.Range("D5").ClearContents()  # RM 14: now missing
$ws.Range("D21").Value = -14.3   # RM 135: was missing, now has a value
$ws.Range("D23").ClearContents() # RM 140: now missing
$ws.Range("D32").Value = -14.7   # SC 193: was missing, now has a value
